$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A1:E31")
$range.Borders.Weight = 4
$range.Borders.Color = 255
